# Auto update stock data
# Update the "as of" date in cell A2 on every worksheet from 2025/10/26 to 2025/10/27.
# The date is stored as literal text (not a real Excel date), so the cell is
# forced to Text format before assignment to keep Excel from re-interpreting
# the "yyyy/mm/dd" looking string as an actual date value.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $cell = $ws.Range("A2")
    if ($cell.Value2 -eq "2025/10/26") {
        $cell.NumberFormat = "@"
        $cell.Value = "2025/10/27"
    }
}
